$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect to make edits, then restore protection.
$ws.Unprotect()

# Update the confidential footer note's "as of" date (2021-05-19 -> 2021-05-20).
$ws.Range("A38").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-20 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) columns for each holding row.
$ws.Range("D2").Value = 0.03568443703457255
$ws.Range("E2").Value = 0.0003317535545022121
$ws.Range("D3").Value = 0.02034674018593421
$ws.Range("E3").Value = 0.002750491159135615
$ws.Range("D4").Value = 0.01942774934160856
$ws.Range("E4").Value = 0.002405773857257376
$ws.Range("D5").Value = 0.03774873277455902
$ws.Range("E5").Value = 0.01200141193081561
$ws.Range("D6").Value = 0.03445093320722065
$ws.Range("E6").Value = 0
$ws.Range("D7").Value = 0.01978710489594781
$ws.Range("E7").Value = 0.002140637140637391
$ws.Range("D8").Value = 0.03708839388156682
$ws.Range("E8").Value = 0.01177982437352765
$ws.Range("D9").Value = 0.02041254345629525
$ws.Range("E9").Value = 0.006959508315256713
$ws.Range("D10").Value = 0.025921568651568
$ws.Range("E10").Value = 0.01097271648873077
$ws.Range("D11").Value = 0.02397155740761354
$ws.Range("E11").Value = 0.006948156066274747
$ws.Range("D12").Value = 0.05727867193010831
$ws.Range("E12").Value = 0.006927854753941753
$ws.Range("D13").Value = 0.02479614823324369
$ws.Range("E13").Value = 0.01488095238095255
$ws.Range("D14").Value = 0.02688688827814161
$ws.Range("E14").Value = 0.001250390747108643
$ws.Range("D15").Value = 0.03267122123694979
$ws.Range("E15").Value = -0.002283907238229133
$ws.Range("D16").Value = 0.01979996830783303
$ws.Range("E16").Value = 0.007280944012051149
$ws.Range("D17").Value = 0.03120530476856416
$ws.Range("E17").Value = 0.01059287239283946
$ws.Range("D18").Value = 0.04199673361585692
$ws.Range("E18").Value = 0.006028286575469499
$ws.Range("D19").Value = 0.1261412818764059
$ws.Range("E19").Value = 0.005336891260840471
$ws.Range("D20").Value = 0.009354314121853308
$ws.Range("E20").Value = 0.003780241935483764
$ws.Range("D21").Value = 0.01534692160616558
$ws.Range("E21").Value = 0.007713884992987419
$ws.Range("D22").Value = 0.01713442337168167
$ws.Range("E22").Value = 0.01632175725835161
$ws.Range("D23").Value = 0.01542881695354638
$ws.Range("E23").Value = 0.004344677769732019
$ws.Range("D24").Value = 0.02143751651487841
$ws.Range("E24").Value = 0.0123355263157896
$ws.Range("D25").Value = 0.01249508781661674
$ws.Range("E25").Value = 0.01621324539708691
$ws.Range("D26").Value = 0.04234394323945196
$ws.Range("E26").Value = 0.00800609988562706
$ws.Range("D27").Value = 0.02404115307829134
$ws.Range("E27").Value = -0.0001961168856638995
$ws.Range("D28").Value = 0.04551956601820722
$ws.Range("E28").Value = 0.007194244604316502
$ws.Range("D29").Value = 0.05543689784114025
$ws.Range("E29").Value = 0.01395180286283737
$ws.Range("D30").Value = 0.01299911831817501
$ws.Range("E30").Value = -0.0006544502617802372
$ws.Range("D31").Value = 0.02072469900128675
$ws.Range("E31").Value = 0.0003832886163281923
$ws.Range("D32").Value = 0.01333746192481331
$ws.Range("E32").Value = 0.009606147934678289
$ws.Range("D33").Value = 0.0419039223054022
$ws.Range("E33").Value = 0.003103983445421621
$ws.Range("D34").Value = 0.01688017880450012
$ws.Range("E34").Value = 0.01547158583754848
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = 0.006687257350440623

# Restore sheet protection.
$ws.Protect()
